# Rework the "債務" (debt) worksheet (sheet2) so it follows the same
# normalized layout as the other property-type sheets:
#   - Row 1 becomes a genuine header row (species, debtor, owner, total,
#     register_date, register_reason, property_category, category, date,
#     legislator_name, legislator_id, source_file, index)
#   - Row 2 keeps the original record's values (species/debtor/owner/
#     total/register_date/register_reason) and gains the extra
#     property_category/category/date/legislator_name/legislator_id/
#     source_file/index columns (H2:N2)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# ---- Row 1 : headers (B1:N1) ----
$ws.Cells.Item(1,2).Value = "species"
$ws.Cells.Item(1,3).Value = "debtor"
$ws.Cells.Item(1,4).Value = "owner"
$ws.Cells.Item(1,5).Value = "total"
$ws.Cells.Item(1,6).Value = "register_date"
$ws.Cells.Item(1,7).Value = "register_reason"
$ws.Cells.Item(1,8).Value = "property_category"
$ws.Cells.Item(1,9).Value = "category"
$ws.Cells.Item(1,10).Value = "date"
$ws.Cells.Item(1,11).Value = "legislator_name"
$ws.Cells.Item(1,12).Value = "legislator_id"
$ws.Cells.Item(1,13).Value = "source_file"
$ws.Cells.Item(1,14).Value = "index"

# ---- Row 2 : data (A2:N2) ----
$ws.Cells.Item(2,1).Value = 81
$ws.Cells.Item(2,2).Value = "中期放款"
$ws.Cells.Item(2,3).Value = "陳亭妃"
$ws.Cells.Item(2,4).Value = "合作金庫商業銀行臺南市北區西門路"
$ws.Cells.Item(2,5).Value = 970000
$ws.Cells.Item(2,6).Value = "89年03月29日"
$ws.Cells.Item(2,7).Value = "信用貸款"
$ws.Cells.Item(2,8).Value = "debt"
$ws.Cells.Item(2,9).Value = "normal"

# column J holds a text date string ("2011-12-28"); pre-format as Text so
# Excel keeps it as a string instead of auto-converting to a date serial
$ws.Range("J2").NumberFormat = "@"
$ws.Cells.Item(2,10).Value = "2011-12-28"

$ws.Cells.Item(2,11).Value = "陳亭妃"
$ws.Cells.Item(2,12).Value = 1708
$ws.Cells.Item(2,13).Value = "tmp1fdf1"
$ws.Cells.Item(2,14).Value = 81

# ---- Formatting : the newly added header cells (H1:N1) get the same
#      bold/centered/bordered style already used for B1:G1 in this sheet
#      (and the header row of the other sheet in this workbook) ----
$headerRange = $ws.Range("H1:N1")
$headerRange.Borders.LineStyle = 1
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
